$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18 ("4.Recycling, Extraction of material
# recycled for each component" - the new "metrec" lookup row), shifting
# the existing GDP projections / AIC / Projections / Merged FD rows down
# by one (old rows 18-21 become 19-22).
$ws.Rows("18:18").Insert()

# Populate the new row with the path-lookup name and Matteo's OneDrive path.
$ws.Range("A18").Value = "metrec"
$ws.Range("E18").Value = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Recycling\Met_rec_comp"

# Restore the selection on the frozen (bottom-right) pane to match the saved file.
$ws.Range("I28").Select()
